$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Nor Shahirah Shaik Amir"
$ws.Range("C2").Value = "N/A"
$ws.Range("D2").Value = "shahirahshaik25@gmail.com"
$ws.Range("G2").Value = "['Malaysia', 'Selangor', 'Sepang']"
$ws.Range("H2").Value = "[['Electric, Electronics & System Engineering (Image processing)', 'Master of Science', 'N/A', 'Universiti Kebangsaan Malaysia', '2015-2019'], ['Electric, Electronics & System Engineering (Communication & Computer)', ""Bachelor's degree"", 'N/A', 'Universiti Kebangsaan Malaysia', '2009-2013'], ['Science majoring in Biology', 'Matriculation', 'N/A', 'Penang Matriculation College', '2008-2009']]"
$ws.Range("I2").Value = "['Google Data Analytics Professional Certificate']"
$ws.Range("J2").Value = "['SQL', 'R (Programming Language)', 'Data Analytics', 'Project Planning', 'Google Sheets', 'Data Analysis', 'Research and Development (R&D)', 'Professional Communication', 'Program Management', 'C (Programming Language)']"
$ws.Range("K2").Value = "['MATLAB', 'Python', 'R', 'SQL', 'Microsoft Word', 'Microsoft PowerPoint', 'Microsoft Excel', 'Microsoft Access', 'Microsoft Publisher', 'Google Docs', 'Google Sheets', 'Google Slides', 'Google Form', 'MPLAB', 'Arduino IDE', 'Visual Studio Code IDE', 'MIT App Inventor', 'Filmora 9', 'SPSS']"
$ws.Range("L2").Value = "['N/A']"
$ws.Range("M2").Value = "[{'job_title': 'Media Officer', 'job_company': 'Pertubuhan IKRAM Malaysia', 'Industries': 'Media', 'start_date': '2021-04', 'end_date': 'Present', 'job_location': 'N/A', 'job_duration': 1.92}, {'job_title': 'Research Officer', 'job_company': 'UKM', 'Industries': 'Research', 'start_date': '2019-04', 'end_date': '2020-08', 'job_location': 'N/A', 'job_duration': 1.42}, {'job_title': 'Research Assistant', 'job_company': 'UKM', 'Industries': 'Research', 'start_date': '2014-04', 'end_date': '2015-07', 'job_location': 'N/A', 'job_duration': 1.33}, {'job_title': 'Internship', 'job_company': 'Astana Digital Sdn. Bhd', 'Industries': 'N/A', 'start_date': '2012-06', 'end_date': '2012-09', 'job_location': 'N/A', 'job_duration': 0.33}]"
$ws.Range("B3").Value = "ELNI FATIENY BINTI MOHD SALMI"
$ws.Range("C3").Value = "011-13244806"
$ws.Range("D3").Value = "elnifatieny136@gmail.com"
$ws.Range("G3").Value = "['Malaysia', 'N/A', 'N/A']"
$ws.Range("H3").Value = "[{'field_of_study': 'Information System', 'level': 'Bachelor of Science', 'cgpa': '3.53', 'university': 'Universiti Teknologi Petronas', 'start_date': '2018-05', 'year_of_graduation': '2021'}]"
$ws.Range("I3").Value = "['N/A']"
$ws.Range("J3").Value = "['Data Analyst', 'Corporate Management', 'Software Designing', 'Project Management', 'Marketing Operation']"
$ws.Range("K3").Value = "['Eclipse', 'Anaconda', 'R studio', 'Android Studio', 'PowerBI', 'Firebase', 'Visual Studio Code', 'Microsoft Office', 'Google Collab', 'Oracle Live', 'phpMyAdmin', 'MySQL', 'SQL Server Management Studio', 'WinSCP', 'Salesforce', 'WordPress']"
$ws.Range("L3").Value = "['Malay', 'English']"
$ws.Range("M3").Value = "[{'job_title': 'Marketing Associate Operation', 'job_company': '2X Marketing Sdn Bhd', 'Industries': 'Marketing', 'start_date': '2022-09', 'end_date': '2024-05-15'}, {'job_title': 'Software Engineer', 'job_company': 'ManagePay System Berhad', 'Industries': 'Finance', 'start_date': '2022-01', 'end_date': '2024-05-15'}, {'job_title': 'Change Management Intern', 'job_company': 'Petroliam Nasional Berhad', 'Industries': 'Oil & Gas', 'start_date': '2020-01', 'end_date': '2020-09'}]"
